# BITACORA.xlsx update (jdre-28/05/15)
# - Adds two new activity log entries (rows 12 & 13) on the "Actividades" sheet
# - Makes "Actividades" the active/selected sheet (was "Datos")
# - Leaves the "Datos" sheet selection at B8

$wb = $excel.ActiveWorkbook

$wsDatos = $wb.Worksheets.Item("Datos")
$wsAct   = $wb.Worksheets.Item("Actividades")

# --- New activity: week 4, 25/05/2015 ---------------------------------
$wsAct.Range("B12").Value = 4
$wsAct.Range("C12").Value = (Get-Date -Year 2015 -Month 5 -Day 25 -Hour 0 -Minute 0 -Second 0)
$wsAct.Range("D12").Value = "Reunion de equipo para opinar sobre el nombre de la compañía, inicio de desarrollo del sitio web empresarial"
$wsAct.Range("E12").Value = "PROCESO"
$wsAct.Range("F12").Value = "Se decidio que el nombre de la compañía sera DevUniverse y se continuara con el desarrollo del sitio web empresarial"
$wsAct.Range("G12").Value = "Sitio WEB fase uno y bitacora"
$wsAct.Range("H12").Value = "Todos"

# D12 gains word-wrap (matches the D10/D11 formatting already used above it)
$wsAct.Range("D12").WrapText = $true

# --- New activity: week 4, 28/05/2015 ----------------------------------
$wsAct.Range("B13").Value = 4
$wsAct.Range("C13").Value = (Get-Date -Year 2015 -Month 5 -Day 28 -Hour 0 -Minute 0 -Second 0)
$wsAct.Range("D13").Value = "Finalizacion de la primera etapa del sitio WEB empresarial"
$wsAct.Range("E13").Value = "FINALIZADA"
$wsAct.Range("F13").Value = "Sitio WEB empresarial DevUniverse en espera de ser evaluado por el profesor lider Alejandro Lara San Juan"
$wsAct.Range("G13").Value = "Sitio WEB empresarial"
$wsAct.Range("H13").Value = "Todos"

# "Responsable" column (H9:H13) is re-formatted so all five rows match:
# centered + word-wrapped
$respRange = $wsAct.Range("H9:H13")
$respRange.HorizontalAlignment = -4108
$respRange.WrapText = $true

# --- View/selection state ----------------------------------------------
# Keep Datos' own selection at B8 ...
$wsDatos.Activate()
$wsDatos.Range("B8").Select()

# ... but the workbook now opens on Actividades, scrolled/selected at H9
$wsAct.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$wsAct.Range("H9").Select()

Write-Output "BITACORA updated"
